# Applies the cell-value changes for Sheets/Siren_Profits.xlsx
# (per-row currentAveragePrice / Leve profit recompute across all 8 job sheets)
$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H62").Value = 125002930
$ws.Range("I62").Value = 142860320
$ws.Range("J62").Value = 1099
$ws.Range("K62").Value = 142860320
$ws.Range("L62").Value = 1099
$ws.Range("M62").Value = -142859696
$ws.Range("N62").Value = -2347
$ws.Range("H65").Value = 125002930
$ws.Range("I65").Value = 142860320
$ws.Range("J65").Value = 1099
$ws.Range("K65").Value = 714301600
$ws.Range("L65").Value = 5495
$ws.Range("M65").Value = -714298480
$ws.Range("N65").Value = -11735
$ws.Range("H86").Value = 47643544
$ws.Range("I86").Value = 3348
$ws.Range("J86").Value = 71463650
$ws.Range("K86").Value = 3348
$ws.Range("L86").Value = 71463650
$ws.Range("M86").Value = -2225
$ws.Range("N86").Value = -71465896
$ws.Range("H89").Value = 47643544
$ws.Range("I89").Value = 3348
$ws.Range("J89").Value = 71463650
$ws.Range("K89").Value = 16740
$ws.Range("L89").Value = 357318250
$ws.Range("M89").Value = -11124
$ws.Range("N89").Value = -357329482
$ws.Range("H106").Value = 9345.091
$ws.Range("I106").Value = 11293
$ws.Range("J106").Value = 5936.25
$ws.Range("K106").Value = 11293
$ws.Range("L106").Value = 5936.25
$ws.Range("M106").Value = -10662
$ws.Range("N106").Value = -7198.25
$ws.Range("H137").Value = 9142.275
$ws.Range("I137").Value = 16479.285
$ws.Range("K137").Value = 49437.855
$ws.Range("M137").Value = -46887.855
$ws.Range("H138").Value = 1452.7028
$ws.Range("I138").Value = 833.0769
$ws.Range("J138").Value = 2917.2727
$ws.Range("K138").Value = 2499.2307
$ws.Range("L138").Value = 8751.8181
$ws.Range("M138").Value = 2640.7693
$ws.Range("N138").Value = -19031.8181

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 62468.707
$ws.Range("J2").Value = 202524.6
$ws.Range("L2").Value = 202524.6
$ws.Range("N2").Value = -202750.6
$ws.Range("H32").Value = 10177.333
$ws.Range("I32").Value = 10192.111
$ws.Range("K32").Value = 10192.111
$ws.Range("M32").Value = -9905.111000000001
$ws.Range("H45").Value = 55290.973
$ws.Range("I45").Value = 115147.836
$ws.Range("K45").Value = 115147.836
$ws.Range("M45").Value = -114770.836
$ws.Range("H61").Value = 9112.857
$ws.Range("I61").Value = 14277.429
$ws.Range("J61").Value = 3948.2856
$ws.Range("K61").Value = 14277.429
$ws.Range("L61").Value = 3948.2856
$ws.Range("M61").Value = -14065.429
$ws.Range("N61").Value = -4372.2856
$ws.Range("H97").Value = 14293690
$ws.Range("J97").Value = 33334316
$ws.Range("L97").Value = 33334316
$ws.Range("N97").Value = -33335308
$ws.Range("H110").Value = 2841.6365
$ws.Range("I110").Value = 1822.5714
$ws.Range("K110").Value = 1822.5714
$ws.Range("M110").Value = 222.4286
$ws.Range("H116").Value = 62468.707
$ws.Range("J116").Value = 202524.6
$ws.Range("L116").Value = 202524.6
$ws.Range("N116").Value = -207112.6
$ws.Range("H132").Value = 3582.1853
$ws.Range("I132").Value = 3836.4443
$ws.Range("K132").Value = 11509.3329
$ws.Range("M132").Value = -8979.332900000001
$ws.Range("H136").Value = 9112.857
$ws.Range("I136").Value = 14277.429
$ws.Range("J136").Value = 3948.2856
$ws.Range("K136").Value = 42832.287
$ws.Range("L136").Value = 11844.8568
$ws.Range("M136").Value = -40282.287
$ws.Range("N136").Value = -16944.8568

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 62468.707
$ws.Range("J3").Value = 202524.6
$ws.Range("L3").Value = 202524.6
$ws.Range("N3").Value = -202752.6
$ws.Range("H94").Value = 11446.333
$ws.Range("I94").Value = 15499.223
$ws.Range("J94").Value = 3340.5557
$ws.Range("K94").Value = 15499.223
$ws.Range("L94").Value = 3340.5557
$ws.Range("M94").Value = -15048.223
$ws.Range("N94").Value = -4242.5557
$ws.Range("H107").Value = 2029.9166
$ws.Range("I107").Value = 2236
$ws.Range("K107").Value = 2236
$ws.Range("M107").Value = -316
$ws.Range("H134").Value = 9222
$ws.Range("I134").Value = 11405.77
$ws.Range("K134").Value = 34217.31
$ws.Range("M134").Value = -31682.31

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 35248.75
$ws.Range("I31").Value = 120000
$ws.Range("J31").Value = 6998.3335
$ws.Range("K31").Value = 120000
$ws.Range("L31").Value = 6998.3335
$ws.Range("M31").Value = -119705
$ws.Range("N31").Value = -7588.3335
$ws.Range("H34").Value = 35248.75
$ws.Range("I34").Value = 120000
$ws.Range("J34").Value = 6998.3335
$ws.Range("K34").Value = 120000
$ws.Range("L34").Value = 6998.3335
$ws.Range("M34").Value = -119798
$ws.Range("N34").Value = -7402.3335
$ws.Range("H99").Value = 5053266.5
$ws.Range("I99").Value = 8297366.5
$ws.Range("K99").Value = 8297366.5
$ws.Range("M99").Value = -8295868.5
$ws.Range("H105").Value = 420548.8
$ws.Range("I105").Value = 525467.25
$ws.Range("K105").Value = 525467.25
$ws.Range("M105").Value = -523720.25
$ws.Range("H122").Value = 12439.889
$ws.Range("I122").Value = 11265
$ws.Range("K122").Value = 33795
$ws.Range("M122").Value = -31345
$ws.Range("H126").Value = 5053266.5
$ws.Range("I126").Value = 8297366.5
$ws.Range("K126").Value = 24892099.5
$ws.Range("M126").Value = -24889629.5

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 55215940
$ws.Range("I4").Value = 53505070
$ws.Range("J4").Value = 60104136
$ws.Range("K4").Value = 160515210
$ws.Range("L4").Value = 180312408
$ws.Range("M4").Value = -160515098
$ws.Range("N4").Value = -180312632
$ws.Range("H131").Value = 1917.3838
$ws.Range("J131").Value = 1922.9072
$ws.Range("L131").Value = 5768.721600000001
$ws.Range("N131").Value = -15848.7216

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H52").Value = 0
$ws.Range("J52").Value = 0
$ws.Range("L52").Value = 0
$ws.Range("N52").ClearContents()
$ws.Range("H102").Value = 8631.182000000001
$ws.Range("I102").Value = 10494.333
$ws.Range("K102").Value = 10494.333
$ws.Range("M102").Value = -8872.333000000001
$ws.Range("H123").Value = 45500
$ws.Range("J123").Value = 45500
$ws.Range("L123").Value = 45500
$ws.Range("N123").Value = -50400
$ws.Range("H126").Value = 8120.174
$ws.Range("I126").Value = 15504.777
$ws.Range("K126").Value = 46514.331
$ws.Range("M126").Value = -44044.331
$ws.Range("H135").Value = 89555
$ws.Range("J135").Value = 89555
$ws.Range("L135").Value = 89555
$ws.Range("N135").Value = -99695

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 24690.818
$ws.Range("I40").Value = 31479.934
$ws.Range("K40").Value = 31479.934
$ws.Range("M40").Value = -31343.934
$ws.Range("H136").Value = 5450.375
$ws.Range("I136").Value = 3844.25
$ws.Range("J136").Value = 7056.5
$ws.Range("K136").Value = 11532.75
$ws.Range("L136").Value = 21169.5
$ws.Range("M136").Value = -8982.75
$ws.Range("N136").Value = -26269.5

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 19346.234
$ws.Range("I107").Value = 1949.2142
$ws.Range("K107").Value = 5847.642599999999
$ws.Range("M107").Value = -3927.642599999999
$ws.Range("H122").Value = 4261.617
$ws.Range("I122").Value = 2129.0303
$ws.Range("K122").Value = 6387.090899999999
$ws.Range("M122").Value = -3937.090899999999
$ws.Range("H126").Value = 33432.848
$ws.Range("I126").Value = 51842.625
$ws.Range("J126").Value = 3977.2
$ws.Range("K126").Value = 155527.875
$ws.Range("L126").Value = 11931.6
$ws.Range("M126").Value = -153057.875
$ws.Range("N126").Value = -16871.6
$ws.Range("H132").Value = 11195.697
$ws.Range("J132").Value = 3324.4375
$ws.Range("L132").Value = 9973.3125
$ws.Range("N132").Value = -15033.3125

